$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.795.03"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").Value = "2.104.28"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.12"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.62"
$ws.Range("E7").Value = "  +2.13%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.388"
$ws.Range("E9").Value = "  +1.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("E12").Value = "  +7.19%  "

$ws.Range("D13").Value = "2.417.19"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.24"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.807"
$ws.Range("E15").Value = "  +3.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "2.125.06"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "38.832.12"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.05"
$ws.Range("E19").Value = "  +2.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.05"
$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.09"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -2.63%  "

$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.10"
$ws.Range("E26").Value = "  +1.30%  "

$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  +5.96%  "

$ws.Range("E29").Value = "  +3.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.36"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("E31").Value = "  +3.82%  "

$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("E33").Value = "  +2.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0621"
$ws.Range("E35").Value = "  +2.25%  "

$ws.Range("E36").Value = "  +2.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("E37").Value = "  +1.01%  "

$ws.Range("E38").Value = "  +1.07%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.34"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("E41").Value = "  +4.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.22"
$ws.Range("E42").Value = "  +1.97%  "

$ws.Range("D43").Value = "1.533.42"
$ws.Range("E43").Value = "  -1.13%  "

$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.79"
$ws.Range("E45").Value = "  +4.52%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("E47").Value = "  +1.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.13"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.04"
$ws.Range("E49").Value = "  +1.38%  "

$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("D51").Value = "2.300.76"
$ws.Range("E51").Value = "  +0.01%  "
